$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 195, pushing the existing rows 195-197
# down to 196-198 (values preserved automatically by the insert).
$ws.Rows("195:195").Insert()

# Populate the newly inserted row 195 with the new record.
$ws.Range("A195").Value = 9
$ws.Range("B195").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C195").Value = "Metropolitana"
$ws.Range("D195").Value = 44656
$ws.Range("E195").Value = 13
$ws.Range("F195").Value = "Fruta"
$ws.Range("G195").Value = 100101
$ws.Range("H195").Value = "Berries"
$ws.Range("I195").Value = 100101001
$ws.Range("J195").Value = "Arándano (blue)"
$ws.Range("K195").Value = "Sin especificar"
$ws.Range("L195").Value = "Primera"
$ws.Range("M195").Value = 280
$ws.Range("N195").Value = 4000
$ws.Range("O195").Value = 4000
$ws.Range("P195").Value = 4000
$ws.Range("Q195").Value = "$/bandeja 2 kilos"
$ws.Range("R195").Value = "Provincia de Linares"
$ws.Range("S195").Value = 2000
$ws.Range("T195").Value = 2

# Give the new date cell the same date style used by the rest of column D.
$ws.Range("D195").NumberFormat = $ws.Range("D196").NumberFormat
